$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (number of interested people) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3354
$ws1.Range("F4").Value = 63
$ws1.Range("F5").Value = 1465
$ws1.Range("F6").Value = 36
$ws1.Range("F7").Value = 322

# Sheet "全部类型" (All types) - same underlying rows duplicated, update column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3354
$ws4.Range("F4").Value = 63
$ws4.Range("F5").Value = 1465
$ws4.Range("F6").Value = 36
$ws4.Range("F8").Value = 322
